$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows above the existing row 233, pushing current 233-236 down to 235-238.
$ws.Range("A233:R234").EntireRow.Insert()

# --- New row 233: Cuatro cascos rojo, Region del Maule ---
$ws.Range("A233").Value = 7
$ws.Range("B233").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C233").Value = "Ñuble"
$ws.Range("D233").Value = 44628
$ws.Range("E233").Value = 16
$ws.Range("F233").Value = 100112002
$ws.Range("G233").Value = "Pimiento"
$ws.Range("H233").Value = "Cuatro cascos rojo"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 120
$ws.Range("K233").Value = 15500
$ws.Range("L233").Value = 16000
$ws.Range("M233").Value = 15750
$ws.Range("N233").Value = "$/caja 15 kilos"
$ws.Range("O233").Value = "Región del Maule"
$ws.Range("P233").Value = 1050
$ws.Range("Q233").Value = 15
$ws.Range("R233").Value = "Hortaliza"

# --- New row 234: Cuatro cascos verde, Region del Maule ---
$ws.Range("A234").Value = 7
$ws.Range("B234").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C234").Value = "Ñuble"
$ws.Range("D234").Value = 44628
$ws.Range("E234").Value = 16
$ws.Range("F234").Value = 100112002
$ws.Range("G234").Value = "Pimiento"
$ws.Range("H234").Value = "Cuatro cascos verde"
$ws.Range("I234").Value = "Primera"
$ws.Range("J234").Value = 120
$ws.Range("K234").Value = 8500
$ws.Range("L234").Value = 9000
$ws.Range("M234").Value = 8750
$ws.Range("N234").Value = "$/caja 15 kilos"
$ws.Range("O234").Value = "Región del Maule"
$ws.Range("P234").Value = 583
$ws.Range("Q234").Value = 15
$ws.Range("R234").Value = "Hortaliza"
